$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay as text (matching the
# original inline-string cell type) instead of being auto-converted to numbers.
$priceRows = @(2,3,5,6,8,11,12,13,16,17,18,19,20,21,22,23,24,27,29,30,31,32,34,36,37,38,39,40,43,44,45,48,50,51)
foreach ($r in $priceRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# Apply the updated price (column D) and volume-change (column E) values.
$ws.Range("D2").Value = "61.886.61"
$ws.Range("E2").Value = "  +3.04%  "
$ws.Range("D3").Value = "3.414.49"
$ws.Range("E3").Value = "  +3.44%  "
$ws.Range("D5").Value = "577.76"
$ws.Range("E5").Value = "  +2.88%  "
$ws.Range("D6").Value = "138.70"
$ws.Range("E6").Value = "  +7.96%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.414.08"
$ws.Range("E8").Value = "  +3.32%  "
$ws.Range("E9").Value = "  +1.59%  "
$ws.Range("E10").Value = "  +10.27%  "
$ws.Range("D11").Value = "7.48"
$ws.Range("E11").Value = "  +1.47%  "
$ws.Range("D12").Value = "0.396"
$ws.Range("E12").Value = "  +7.04%  "
$ws.Range("D13").Value = "3.998.05"
$ws.Range("E13").Value = "  +3.68%  "
$ws.Range("E14").Value = "  +2.15%  "
$ws.Range("E15").Value = "  +9.26%  "
$ws.Range("D16").Value = "3.413.45"
$ws.Range("E16").Value = "  +3.57%  "
$ws.Range("D17").Value = "25.48"
$ws.Range("E17").Value = "  +6.35%  "
$ws.Range("D18").Value = "61.941.63"
$ws.Range("E18").Value = "  +2.77%  "
$ws.Range("D19").Value = "14.12"
$ws.Range("E19").Value = "  +7.00%  "
$ws.Range("D20").Value = "5.91"
$ws.Range("E20").Value = "  +5.78%  "
$ws.Range("D21").Value = "9.51"
$ws.Range("E21").Value = "  +7.80%  "
$ws.Range("D22").Value = "390.84"
$ws.Range("E22").Value = "  +12.21%  "
$ws.Range("D23").Value = "0.574"
$ws.Range("E23").Value = "  +4.15%  "
$ws.Range("D24").Value = "3.550.81"
$ws.Range("E24").Value = "  +3.57%  "
$ws.Range("E25").Value = "  +19.49%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "71.68"
$ws.Range("E27").Value = "  +4.53%  "
$ws.Range("E28").Value = "  +10.76%  "
$ws.Range("D29").Value = "7.67"
$ws.Range("E29").Value = "  +5.37%  "
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("D31").Value = "8.31"
$ws.Range("E31").Value = "  +6.68%  "
$ws.Range("D32").Value = "0.160"
$ws.Range("E32").Value = "  +5.72%  "
$ws.Range("E33").Value = "  +4.31%  "
$ws.Range("D34").Value = "3.446.21"
$ws.Range("E34").Value = "  +3.57%  "
$ws.Range("D36").Value = "23.61"
$ws.Range("E36").Value = "  +4.61%  "
$ws.Range("D37").Value = "5.50"
$ws.Range("E37").Value = "  +4.28%  "
$ws.Range("D38").Value = "7.01"
$ws.Range("E38").Value = "  +4.28%  "
$ws.Range("D39").Value = "1.56"
$ws.Range("E39").Value = "  +6.39%  "
$ws.Range("D40").Value = "162.68"
$ws.Range("E40").Value = "  +4.03%  "
$ws.Range("E41").Value = "  +6.43%  "
$ws.Range("E42").Value = "  +15.60%  "
$ws.Range("D43").Value = "0.793"
$ws.Range("E43").Value = "  +7.22%  "
$ws.Range("D44").Value = "25.55"
$ws.Range("E44").Value = "  +13.63%  "
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("E46").Value = "  +6.04%  "
$ws.Range("E47").Value = "  +4.75%  "
$ws.Range("D48").Value = "41.63"
$ws.Range("E48").Value = "  +2.40%  "
$ws.Range("E49").Value = "  +4.41%  "
$ws.Range("D50").Value = "23.13"
$ws.Range("E50").Value = "  +6.06%  "
$ws.Range("D51").Value = "2.398.37"
$ws.Range("E51").Value = "  +11.58%  "
